$d = $word.ActiveDocument

# The title run currently reads "eCoacing Log System" (misspelled). The
# author's fix splits it into three runs with identical formatting:
#   "eCoac" + "h" + "ing Log System"
# i.e. an "h" is inserted between "eCoac" and "ing Log System" so the
# word reads "eCoaching".

$rng = $d.Content
$rng.Find.Execute("eCoac", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($rng.Find.Found) {
    $insertPoint = $d.Range($rng.End, $rng.End)
    $insertPoint.InsertAfter("h")

    # Re-apply (toggle off/on) direct character formatting that already
    # matches the surrounding run. This does not change the visible
    # formatting (Bold was already on), but it forces the newly typed
    # "h" to stay in its own run instead of being silently re-merged
    # with its neighbours, matching the run layout in the target
    # document (three runs: "eCoac", "h", "ing Log System").
    $hRange = $d.Range($rng.End, $rng.End + 1)
    $hRange.Bold = 0
    $hRange.Bold = 1
}
